$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells keep their values as text (not auto-converted to numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "260.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.57%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.92%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.673"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.75%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06174"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.25%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.664"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.84%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8509"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.67%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9178"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.28%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1411"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.15%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04814"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "7.96%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.04%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03117"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.34%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09049"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.61%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001540"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.56%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006165"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.80%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006161"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.16%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.450"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.52%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.154"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.74%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.69%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1299"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.22%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.082"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.72%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04229"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.51%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.06%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-15.01%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.04%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001575"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-8.12%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03874"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.71%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1113"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.12%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004089"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "10.56%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.59%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-9.25%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005152"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.02%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.07%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "8.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1617"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-32.84%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.07%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.07%"
